# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) values scraped for
# 苏州-漫展信息.xlsx across the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = 89
$ws1.Range("F9").Value = 1179
$ws1.Range("F10").Value = 16647
$ws1.Range("F14").Value = 6451
$ws1.Range("F15").Value = 649
$ws1.Range("F21").Value = 63
$ws1.Range("F23").Value = 643
$ws1.Range("F28").Value = 235
$ws1.Range("F29").Value = 910
$ws1.Range("F30").Value = 68
$ws1.Range("F31").Value = 5075
$ws1.Range("F33").Value = 11452
$ws1.Range("F38").Value = 3857

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 89
$ws4.Range("F9").Value = 1179
$ws4.Range("F10").Value = 16647
$ws4.Range("F14").Value = 6451
$ws4.Range("F15").Value = 649
$ws4.Range("F21").Value = 63
$ws4.Range("F23").Value = 643
$ws4.Range("F28").Value = 235
$ws4.Range("F29").Value = 910
$ws4.Range("F30").Value = 68
$ws4.Range("F31").Value = 5075
$ws4.Range("F34").Value = 11452
$ws4.Range("F39").Value = 3857

$wb.Save()
